# The KeywordTok/ImportTok/CommentTok/... character styles in styles.xml had
# their <w:rPr> children in the wrong order relative to wml.xsd (CT_RPr):
# <w:color> was written before <w:b/>/<w:i/>, but the schema's declared
# sequence requires b/bCs/i/iCs to precede color. OOXMLValidatorCLI flags
# this as Sch_UnexpectedElementContentExpectingComplex even though xmllint
# stays quiet. Re-assigning the (already-correct) Bold/Italic values forces
# the style's rPr to be rewritten in schema order, fixing the element
# sequence without changing any actual formatting.

$d = $word.ActiveDocument

# styleId -> which Font toggles to re-apply (bold / italic), in the order
# they must appear before <w:color> per the schema.
$fixes = [ordered]@{
    "KeywordTok"       = @("b")
    "ImportTok"        = @("b")
    "CommentTok"       = @("i")
    "DocumentationTok" = @("i")
    "AnnotationTok"    = @("b", "i")
    "CommentVarTok"    = @("b", "i")
    "ControlFlowTok"   = @("b")
    "InformationTok"   = @("b", "i")
    "WarningTok"       = @("b", "i")
    "AlertTok"         = @("b")
    "ErrorTok"         = @("b")
}

foreach ($styleId in $fixes.Keys) {
    $style = $d.Styles($styleId)
    foreach ($toggle in $fixes[$styleId]) {
        if ($toggle -eq "b") {
            $style.Font.Bold = $style.Font.Bold
        } elseif ($toggle -eq "i") {
            $style.Font.Italic = $style.Font.Italic
        }
    }
}
